$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: label changes to "unanswerable:other"; remark gets a numeric-looking
# value that must stay text, so force Text number format before assigning.
$ws.Range("R2").Value = 'unanswerable:other'
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = '2222222222'

# Row 3: label -> answerable; remark becomes a real number 233
$ws.Range("R3").Value = 'answerable'
$ws.Range("T3").Value = 233

# Row 4: answer_consistency 0 -> 1
$ws.Range("S4").Value = 1

# Row 5: answer_consistency 1 -> 0
$ws.Range("S5").Value = 0

# Row 6: answer_consistency 0 -> 1
$ws.Range("S6").Value = 1

# Row 7: label -> answerable
$ws.Range("R7").Value = 'answerable'

# Row 8: label -> answerable; answer_consistency -> 1; remark -> big number
$ws.Range("R8").Value = 'answerable'
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 1111111111111

# Row 9: label -> unanswerable:understandability
$ws.Range("R9").Value = 'unanswerable:understandability'

# Row 11: label -> unanswerable:other; remark text (keep trailing ".0" as text)
$ws.Range("R11").Value = 'unanswerable:other'
$ws.Range("T11").NumberFormat = "@"
$ws.Range("T11").Value = '1232132123.0'

# Row 13: remark switches from a number to text "213213233333.0"
$ws.Range("T13").NumberFormat = "@"
$ws.Range("T13").Value = '213213233333.0'

# Row 14: label -> unanswerable:information sufficiency
$ws.Range("R14").Value = 'unanswerable:information sufficiency'

# Row 15: label -> unanswerable:factual consistency
$ws.Range("R15").Value = 'unanswerable:factual consistency'

# Row 16: fill in previously-empty label/answer_consistency
$ws.Range("R16").Value = 'unanswerable:contextual relevance'
$ws.Range("S16").Value = 0

# Row 17
$ws.Range("R17").Value = 'unanswerable:specificity'
$ws.Range("S17").Value = 0

# Row 18
$ws.Range("R18").Value = 'unanswerable:understandability'
$ws.Range("S18").Value = 0

# Row 19 (also gets a Chinese-text remark)
$ws.Range("R19").Value = 'unanswerable:other'
$ws.Range("S19").Value = 0
$ws.Range("T19").Value = '啊啊啊啊啊啊啊啊啊啊啊啊'

# Row 20
$ws.Range("R20").Value = 'unanswerable:understandability'
$ws.Range("S20").Value = 0

# Rows 21-101: label -> answerable, answer_consistency -> 1
for ($r = 21; $r -le 101; $r++) {
    $ws.Range("R$r").Value = 'answerable'
    $ws.Range("S$r").Value = 1
}
